$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 546-547 (existing rows 546-566 shift down to 548-568),
# carrying formatting (incl. the date style on column D) along with them.
$ws.Rows("546:547").Insert()

# Populate the two newly inserted rows with the new weekly price records.
$ws.Range("A546").Value = 4
$ws.Range("B546").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C546").Value = "Los Lagos"
$ws.Range("D546").Value = 45041
$ws.Range("E546").Value = 10
$ws.Range("F546").Value = "Fruta"
$ws.Range("G546").Value = 100102
$ws.Range("H546").Value = "Cítricos"
$ws.Range("I546").Value = 100102006
$ws.Range("J546").Value = "Pomelo"
$ws.Range("K546").Value = "Start Ruby"
$ws.Range("L546").Value = "Primera"
$ws.Range("M546").Value = 200
$ws.Range("N546").Value = 15000
$ws.Range("O546").Value = 16000
$ws.Range("P546").Value = 15500
$ws.Range("Q546").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R546").Value = "Región de O'Higgins"
$ws.Range("S546").Value = 1107
$ws.Range("T546").Value = 14

$ws.Range("A547").Value = 4
$ws.Range("B547").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C547").Value = "Los Lagos"
$ws.Range("D547").Value = 45041
$ws.Range("E547").Value = 10
$ws.Range("F547").Value = "Fruta"
$ws.Range("G547").Value = 100102
$ws.Range("H547").Value = "Cítricos"
$ws.Range("I547").Value = 100102006
$ws.Range("J547").Value = "Pomelo"
$ws.Range("K547").Value = "Start Ruby"
$ws.Range("L547").Value = "Segunda"
$ws.Range("M547").Value = 100
$ws.Range("N547").Value = 13000
$ws.Range("O547").Value = 13000
$ws.Range("P547").Value = 13000
$ws.Range("Q547").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R547").Value = "Región de O'Higgins"
$ws.Range("S547").Value = 929
$ws.Range("T547").Value = 14
